$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (dotted thousands separators,
# variable decimal places, sometimes with subscript-zero notation).
# Where the new price string is purely numeric-looking (e.g. "219.19"),
# pre-set the cell to Text format so the assignment is not silently
# auto-converted to a Number (which would also strip meaningful
# trailing/leading zeros, e.g. "20.80" -> 20.8). Column E percentage
# strings keep their padding spaces/sign, so Excel already stores them
# as text without any extra nudging.

$ws.Range("D2").Value = '26.338.62'
$ws.Range("E2").Value = '  +1.05%  '

$ws.Range("D3").Value = '1.666.19'
$ws.Range("E3").Value = '  +0.90%  '

$ws.Range("E4").Value = '  +1.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.19'
$ws.Range("E5").Value = '  +0.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5338'
$ws.Range("E6").Value = '  +1.33%  '

$ws.Range("E7").Value = '  +0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2663'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06386'
$ws.Range("E9").Value = '  +1.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.80'
$ws.Range("E10").Value = '  +2.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07851'
$ws.Range("E11").Value = '  +0.82%  '

$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("D13").Value = '1.674.68'
$ws.Range("E13").Value = '  +0.34%  '

$ws.Range("D14").Value = '1.894.81'
$ws.Range("E14").Value = '  +0.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5536'
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.88'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '26.360.16'
$ws.Range("E18").Value = '  +1.09%  '

$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.676'
$ws.Range("E20").Value = '  +2.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.54'
$ws.Range("E21").Value = '  +1.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.29'
$ws.Range("E22").Value = '  +2.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.032'
$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.013'
$ws.Range("E24").Value = '  +1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.07'
$ws.Range("E25").Value = '  +2.25%  '

$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.500'
$ws.Range("E29").Value = '  +5.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05850'
$ws.Range("E30").Value = '  +0.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  +0.72%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.587'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.277'
$ws.Range("E33").Value = '  +0.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.603'
$ws.Range("E34").Value = '  +1.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9694'
$ws.Range("E35").Value = '  +2.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.830'
$ws.Range("E36").Value = '  +1.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.418'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5825'
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01606'
$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8612'
$ws.Range("E40").Value = '  +1.53%  '

$ws.Range("D41").Value = '1.064.21'
$ws.Range("E41").Value = '  +3.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.830'
$ws.Range("E42").Value = '  +1.68%  '

$ws.Range("E43").Value = '  +0.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '104.64'
$ws.Range("E44").Value = '  -0.19%  '

$ws.Range("D45").Value = '1.806.25'
$ws.Range("E45").Value = '  +0.71%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.79'
$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.976'
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05166'

# Rows 47-49 reordering: BabyDogeCoin moves to rank 47, Frax to 48, Mantle to 49
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈107'
$ws.Range("E47").Value = '  -4.48%  '

$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.013'
$ws.Range("E48").Value = '  +1.17%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4394'
$ws.Range("E49").Value = '  +1.54%  '
